$wb = $excel.ActiveWorkbook

$contact = $wb.Worksheets.Item("Contact")
$contact.Range("A2").Value = "Activity Test External Contact"
$contact.Range("B2").Value = "ActivityCompany"

$contact.Activate()
$contact.Range("A2:B2").Select()

$wb.Save()
